$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KARINA's (account 005141215, row 7) Saldo value before any row
# deletions shift row numbers.
$ws.Cells.Item(7, 3).Value = 8512.18

# Delete rows bottom-to-top so earlier row numbers stay stable while we work:
#   row 13 -> 002064834 / RAFAELA / 4000
#   row 10 -> 004508516 / EDUARDO  / 5000
#   row 9  -> 004508504 / FERNANDO / 5000
#   row 8  -> 004500804 / RAFAEL   / 5000
#   row 3  -> 005206566 / LEVI     / 60362.28
$ws.Rows(13).Delete()
$ws.Rows(10).Delete()
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()
$ws.Rows(3).Delete()
